$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.645.53'
$ws.Range('E2').Value = '  +0.85%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.478.89'
$ws.Range('E3').Value = '  +0.23%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.08'
$ws.Range('E5').Value = '  +0.08%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.78'
$ws.Range('E6').Value = '  +3.50%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.620'
$ws.Range('E7').Value = '  +5.89%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.476.89'
$ws.Range('E9').Value = '  +0.23%  '

$ws.Range('E10').Value = '  +8.97%  '

$ws.Range('E11').Value = '  -1.04%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.429'
$ws.Range('E12').Value = '  +1.56%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.078.02'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.13'
$ws.Range('E14').Value = '  +2.72%  '

$ws.Range('E15').Value = '  -0.94%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.632.93'
$ws.Range('E16').Value = '  +0.90%  '

$ws.Range('E17').Value = '  +1.52%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.479.16'
$ws.Range('E18').Value = '  +0.60%  '

$ws.Range('E19').Value = '  -0.32%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.11'
$ws.Range('E20').Value = '  -1.25%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '395.00'
$ws.Range('E21').Value = '  +2.34%  '

$ws.Range('E22').Value = '  +1.50%  '

$ws.Range('E23').Value = '  +1.86%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.23%  '

$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.539'
$ws.Range('E25').Value = '  +1.49%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '72.02'
$ws.Range('E26').Value = '  -0.72%  '

$ws.Range('E27').Value = '  +0.72%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.36'
$ws.Range('E28').Value = '  +1.12%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.176'
$ws.Range('E29').Value = '  -0.43%  '

$ws.Range('E30').Value = '  +0.38%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.13'
$ws.Range('E31').Value = '  +0.79%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.40'
$ws.Range('E32').Value = '  -0.45%  '

$ws.Range('E33').Value = '  +1.31%  '

$ws.Range('E34').Value = '  +0.73%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.34'
$ws.Range('E35').Value = '  +1.75%  '

$ws.Range('E36').Value = '  -0.09%  '

$ws.Range('E37').Value = '  -1.81%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '161.85'
$ws.Range('E38').Value = '  -1.63%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.892'
$ws.Range('E39').Value = '  +2.80%  '

$ws.Range('E40').Value = '  +11.56%  '

$ws.Range('E41').Value = '  -2.04%  '

$ws.Range('E42').Value = '  +2.36%  '

$ws.Range('E43').Value = '  -2.14%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.31'
$ws.Range('E44').Value = '  +0.99%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0716'
$ws.Range('E45').Value = '  +0.02%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.32'
$ws.Range('E46').Value = '  -2.73%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.746.20'
$ws.Range('E47').Value = '  -1.50%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '41.60'
$ws.Range('E48').Value = '  -1.29%  '

$ws.Range('E49').Value = '  +0.67%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '329.32'
$ws.Range('E50').Value = '  -2.69%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.04'
$ws.Range('E51').Value = '  -1.61%  '
